$wb = $excel.ActiveWorkbook

# --- YDS sheet: weekly log strings (append Week 15 logged + Week 16 simulated) ---
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = '1 6 7 4 8 9 2 -3 3 2 3 1 1 2 4 3 2 -2 4 9 4 8 3 12 1 1 7 3 8 3 4 13 0 0 0 5 2 8 4 0 5 12 1 14 6 2 9 -1 1 1 9 3 3 1 4 2 23 -1 2 2 6 9 -7 5 6 6 1 5 9 5 2 4 8 0 2 10 1 1 23 38 1 6 9 -1 0 12 -5 6 2 4 7 5 0 -5 6 0 -3 -1 1 10 3 11 4 4 3 13 2 2 -1 2 2 19 13 0 4 4 8 22 -1 7 2 1 -1 3 -4 1 1 3 2 3 4 3 3 1 2 0 83 2 -1 4 3 13 1 12 -1 10 11 0 1 8 2 6 2 0 33 -2 4 16 1 1 5 10 -2 0 1 9 8 3 3 0 5 0 8 5 17 -1 0 5 0 5 3 -1 2 7 3 12 11 0 32 2 6 5 4 2 3 1 6 5 4 1 4 4 9 4 34 2 3 12 3 13 21 1 3 6 0 1 9 11 2 5 2 0 3 25 10 78 1 0 0 34 0 8 3 7 20 5 8 4 4 9 3 1 3 6 4 0 6 0 -1 3 2 0 1 3 8 10 6 10 3 3 11 2 4 5 1 6 0 18 4 3 7 -1 2 2 3 5 40 10 9 0 3 2 1 19 6 5 18 9 3 -3 1 5 0 4 -2 0 5 12 5 1 16 3 2 5 15 5 10 15 1 3 4 2 10 8 6 2 1 10 1 -1 6 0 7 6 0 3 8 2 0 7 0 -1 15 11 -1 4 9 3 16 4 1 11 25 4 5 3 6 -3 7 7 3 6 10 0 12 3 37 11 4 6 4 5 3 4 2 3 1 11 7 2 3 1 0 9 7 1 6 3 7 -2 2 12 6 3 0 0 -4 1 2 67'
$ydsWs.Range("C2").Value = '4 2 11 5 4 33 4 5 4 4 2 13 3 2 0 5 9 1 1 4 6 9 5 4 -3 6 -3 2 2 4 13 2 -4 3 -3 5 1 9 10 6 4 2 4 0 5 2 10 5 11 6 2 -5 1 2 2 1 3 4 17 0 19 1 -1 4 2 1 3 4 8 2 8 2 9 5 1 2 3 2 28 10 9 2 2 12 1 5 1 -1 13 0 4 2 4 2 0 6 3 -2 4 -1 4 3 2 0 0 4 0 1 3 3 2 -2 13 4 3 -1 9 2 4 0 12 5 -1 3 6 4 -3 6 6 6 2 3 0 0 5 -1 24 5 1 6 1 0 6 6 3 0 11 2 1 3 1 2 2 0 6 35 2 -1 0 4 5 3 20 14 14 12 6 2 3 6 5 3 -1 2 3 -3 2 0 4 1 5 4 3 1 3 4 1 13 2 5 9 2 1 0 6 1 1 2 2 -1 4 2 4 5 3 3 6 -1 -3 2 4 5 -2 5 3 16 3 2 0 5 3 4 4 4 3 1 11 9 4 5 5 1 2 0 66 3 5 12 4 5 -1 8 8 10 3 8 14 0 18 3 1 1 1 2 3 10 0 3 8 9 4 28 3 8 5 2 4 7 3 7 -3 4 2 2 1 11 6 12 4 1 4 5 8 8 1 2 3 1 11 8 3 28 4 2 4 2 4 0 5 7 0 7 5 5 1 7 16 2 2 8 3 4 12 1 4 -2 4 -3 -2 13 9 3 2 5 3 6 11 0 18 2'
$ydsWs.Range("B3").Value = '10 6 10 7 24 6 10 14 10 8 4 5 7 15 16 8 9 7 8 13 7 1 16 11 19 14 10 9 5 9 42 2 18 7 13 7 3 7 12 23 8 0 15 9 34 10 10 8 8 7 12 4 18 4 3 11 4 5 27 36 18 3 4 9 3 1 7 6 12 1 5 4 28 6 3 2 3 17 8 -2 3 14 10 1 3 41 11 24 20 76 11 15 6 8 15 11 7 -2 5 9 16 -1 31 42 21 24 18 4 4 29 23 2 15 12 9 17 13 7 51 52 6 18 28 10 -1 57 11 -2 2 -1 -1 5 5 -2 2 14 13 9 5 6 28 37 5 8 0 7 10 7 7 7 8 9 7 12 5 6 1 9 6 5 -2 11 3 9 38 8 9 8 6 13 13 9 20 12 1 10 2 17 3 12 11 12 27 31 2 28 19 7 4 7 5 7 12 2 18 15 9 3 6 10 6 15 -5 1 -3 10 16 3 27 12 8 -4 18 18 5 23 9 18 0 3 8 8 4 9 8 62 10 14 3 3 15 9 5 12 22 9 23 -3 15 5 15 11 11 7 8 7 10 12 4 8 9 14 6 13 2 8 17 24 14 4 17 0 6 7 8 23 0 19'
$ydsWs.Range("C3").Value = '23 11 6 22 9 16 4 6 69 4 12 5 3 9 30 6 4 15 16 23 14 13 16 9 4 8 16 11 14 43 2 11 44 10 14 7 3 6 6 9 3 25 7 18 26 14 13 13 -2 10 11 13 7 3 5 2 18 4 10 8 6 4 -7 11 8 7 4 22 25 42 22 4 5 5 17 15 28 7 28 2 6 11 12 8 8 7 17 3 19 7 43 27 18 7 6 15 5 17 10 9 20 8 4 14 4 11 7 7 5 6 3 9 10 14 5 8 -4 18 13 9 4 5 5 3 20 11 17 8 1 9 13 8 6 3 14 7 6 12 9 9 0 28 10 10 16 6 1 29 27 14 3 9 5 5 8 7 3 8 8 6 14 7 57 13 10 11 7 5 10 11 14 9 10 10 13 11 13 12 0 11 12 15 26 19 -5 18 26 7 19 20 10 19 -1 18 12 12 0 13 -3 14 12 5 22 12 10 14 6 8 4 11 0 12 24 15 9 5 6 13 18 0 23 4 14 9 9 9 7 11 9 7 20 2 18 3 6 14 12 7 6 10 6 7 16 7 31 1 9 20 9 0 0 4 -1 10 8 7 7 15 9 26 5 4 11 10 32 9 2 9 3 8 6 13 15 5 12 11 2 10 4 14 2 5 4 10 9 4 3 5 15 25 7 9 10 14 10 8 11 13 12 12 9 20 4 12 43 7 9 16'

# --- OFF sheet ---
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 178
$offWs.Range("F2").Value = 68
$offWs.Range("G2").Value = 55
$offWs.Range("H2").Value = 4
$offWs.Range("J2").Value = 30
$offWs.Range("L2").Value = 265
$offWs.Range("M2").Value = 160
$offWs.Range("O2").Value = 21
$offWs.Range("P2").Value = 14
$offWs.Range("Q2").Value = 492
$offWs.Range("C3").Value = 135
$offWs.Range("F3").Value = 69
$offWs.Range("G3").Value = 17
$offWs.Range("I3").Value = 48
$offWs.Range("J3").Value = 40
$offWs.Range("N3").Value = 15

# --- DEF sheet ---
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 160
$defWs.Range("E2").Value = 6
$defWs.Range("F2").Value = 55
$defWs.Range("G2").Value = 45
$defWs.Range("J2").Value = 22
$defWs.Range("L2").Value = 297
$defWs.Range("M2").Value = 191
$defWs.Range("O2").Value = 15
$defWs.Range("P2").Value = 8
$defWs.Range("Q2").Value = 523
$defWs.Range("B3").Value = 9
$defWs.Range("C3").Value = 170
$defWs.Range("D3").Value = 3
$defWs.Range("E3").Value = 21
$defWs.Range("F3").Value = 101
$defWs.Range("G3").Value = 28
$defWs.Range("H3").Value = 20
$defWs.Range("I3").Value = 57
$defWs.Range("J3").Value = 47
$defWs.Range("N3").Value = 17

# --- ST sheet ---
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 81
$stWs.Range("D2").Value = 48
$stWs.Range("F2").Value = 127
$stWs.Range("G2").Value = 122
$stWs.Range("J2").Value = 41
$stWs.Range("K2").Value = 40
$stWs.Range("L2").Value = 32
$stWs.Range("M2").Value = 23
$stWs.Range("B3").Value = 41
$stWs.Range("B4").Value = '68 64 60 63 62 62 65 62 62 66 60 62 63 63 62 53 63 59 61 60 67 70 63 63 62 65 65 60 62 62 67 58 61 59 67 65 61 68 65 65'
$stWs.Range("B5").Value = '22 21 21 28 15 16 23 23 21 23 27 23 34 15 19 10 22 27 27 33 27 25 21 28 19 20 25 22 12 12 27 20 17 17 24 30 12 26 22 20'
$stWs.Range("B6").Value = '22 19 41 16 17 21 23 0 26 20 27 31 14 18 27 31 25 18 3 15 27 19 26 33 19 72 21 22 20'
$stWs.Range("D3").Value = '46 47 52 42 45 40 40 50 56 51 56 57 37 38 52 42 40 50 49 36 34 79 46 35 42 37 63 46 35 38 47 38 47 43 43 42 43 55 46 31 41 47 46 40 43 41 46 41'
$stWs.Range("D4").Value = '6 0 9 0 0 0 12 18 0 14 -1 14 0 0 7 0 0 3 0 0 0 0 0 0 13 8 0 2 19 0 0 0 0 0 0 8 6 10 12 0 0 -1 4 0 0 4 1 0'
$stWs.Range("D5").Value = '0 7 0 0 0 0 0 0 0 0 0 0 -1 6 0 13 11 6 0 0 0 2 5 0 0 0 11 16 0 0 0 10 0 0 0 0 11 0 0 0 7 0 9 13 7 -1 16 0 0 15 0 0'

# --- TURNS sheet ---
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 8
$turnsWs.Range("C2").Value = 9
$turnsWs.Range("E2").Value = 8
$turnsWs.Range("D3").Value = 6
$turnsWs.Range("E3").Value = 7

# --- PEN sheet ---
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("D2").Value = 3

